# Updated cryptos list on Thu Sep 14 15:15:32 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# swaps two rows whose rankings crossed (ImmutableX/MXToken at 40/41 and
# EnergySwap/Cronos at 49/50). All of these sheet cells are stored as plain
# text (prices like "26.736.07" / "0.490" and padded percentages like
# "  +1.78%  " are not real numbers), so a helper is used to force each
# write to land as text -- otherwise Excel's automatic type detection would
# silently reinterpret single-dot decimals (e.g. "0.491", "4.06") as
# numbers and normalize/round them, which would corrupt the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    # Force text interpretation so things like "0.491", "4.06", "0.0518"
    # aren't silently parsed as numbers (which would drop trailing zeros /
    # introduce float rounding), then restore the default "Normal" style so
    # no stray per-cell number format sticks around afterwards.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-CellText "D2" '26.714.91'
Set-CellText "E2" '  +1.56%  '
Set-CellText "D3" '1.637.22'
Set-CellText "E3" '  +1.70%  '
Set-CellText "E4" '  +0.00%  '
Set-CellText "D5" '213.71'
Set-CellText "E5" '  +0.33%  '
Set-CellText "E6" '  -0.04%  '
Set-CellText "D7" '0.491'
Set-CellText "E7" '  +1.06%  '
Set-CellText "E8" '  +0.70%  '
Set-CellText "E9" '  +0.96%  '
Set-CellText "D10" '19.11'
Set-CellText "E10" '  +4.04%  '
Set-CellText "D11" '0.0833'
Set-CellText "E11" '  +2.28%  '
Set-CellText "D12" '1.864.05'
Set-CellText "E12" '  +1.73%  '
Set-CellText "D13" '1.637.80'
Set-CellText "E13" '  +1.75%  '
Set-CellText "D14" '4.06'
Set-CellText "E14" '  +0.38%  '
Set-CellText "D15" '0.525'
Set-CellText "E15" '  +1.94%  '
Set-CellText "D16" '26.687.58'
Set-CellText "E16" '  +1.49%  '
Set-CellText "D17" '63.36'
Set-CellText "E17" '  +2.53%  '
Set-CellText "D18" '0.0₃0734'
Set-CellText "E18" '  +0.67%  '
Set-CellText "D19" '208.89'
Set-CellText "E19" '  +2.57%  '
Set-CellText "E20" '  -0.01%  '
Set-CellText "E21" '  +0.86%  '
Set-CellText "E22" '  +1.27%  '
Set-CellText "E23" '  +1.57%  '
Set-CellText "E24" '  -1.23%  '
Set-CellText "D25" '145.80'
Set-CellText "E25" '  +1.02%  '
Set-CellText "E26" '  -0.03%  '
Set-CellText "E27" '  -1.45%  '
Set-CellText "D28" '15.40'
Set-CellText "E28" '  +0.93%  '
Set-CellText "E29" '  +1.68%  '
Set-CellText "D30" '0.0518'
Set-CellText "E30" '  +5.95%  '
Set-CellText "E31" '  +0.10%  '
Set-CellText "D32" '3.23'
Set-CellText "E32" '  +0.92%  '
Set-CellText "D33" '2.96'
Set-CellText "E33" '  +0.45%  '
Set-CellText "E34" '  +1.64%  '
Set-CellText "E35" '  -0.31%  '
Set-CellText "D36" '1.166.91'
Set-CellText "E36" '  +0.68%  '
Set-CellText "E37" '  +0.58%  '
Set-CellText "D38" '0.813'
Set-CellText "E38" '  +1.82%  '
Set-CellText "E39" '  -0.01%  '
Set-CellText "B40" 'MXToken'
Set-CellText "C40" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText "D40" '2.33'
Set-CellText "E40" '  -0.03%  '
Set-CellText "B41" 'ImmutableX'
Set-CellText "C41" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText "D41" '0.504'
Set-CellText "E41" '  +0.39%  '
Set-CellText "E42" '  +3.13%  '
Set-CellText "D43" '0.792'
Set-CellText "E43" '  +0.65%  '
Set-CellText "D44" '1.774.11'
Set-CellText "E44" '  +1.66%  '
Set-CellText "D45" '92.57'
Set-CellText "E45" '  +0.91%  '
Set-CellText "D46" '1.55'
Set-CellText "E46" '  +0.89%  '
Set-CellText "D47" '54.74'
Set-CellText "E47" '  +0.65%  '
Set-CellText "D48" '0.0₆0103'
Set-CellText "E48" '  +7.07%  '
Set-CellText "B49" 'Cronos'
Set-CellText "C49" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText "D49" '0.0511'
Set-CellText "E49" '  +0.92%  '
Set-CellText "B50" 'EnergySwap'
Set-CellText "C50" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText "D50" '7.63'
Set-CellText "E50" '  +5.13%  '
Set-CellText "D51" '0.410'
Set-CellText "E51" '  +0.85%  '
